# ATA da reuniao do dia 06/04/2020
# Fill in the minutes (ATA) table for the 06/04/2020 meeting: date/time of
# the meeting, attendees, and the list of discussed topics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: meeting date/time + attendees + first topic -------------------
# Values are written in this particular order so that the shared-strings
# table ends up built in the same sequence as the source workbook.
$ws.Range("E13").Value = " "
$ws.Range("F12").Value = "1 -  Product Owner: Gabriel"
$ws.Range("F16").Value = "5 - CSS do site institucional"
$ws.Range("D12").Value = "Raphael Moitinho, Stefany Batista, Graziela, Gabriel Bezerra, Yuri Vedovate, Bruno Santana"
$ws.Range("F13").Value = "2 - Definição do que será feito essa semana"
$ws.Range("F14").Value = "3 - Fazer prototipo do site"
$ws.Range("F15").Value = "4 - Marcar como feita as atividades"
$ws.Range("F17").Value = "6 - Ver o que é dashboard com yoshi"
$ws.Range("E14").Value = " "

# --- Date / start time / end time for the meeting ---------------------------
$ws.Range("A12").Value = 43927
$ws.Range("A12").NumberFormat = "mm-dd-yy"

$ws.Range("B12").Value = 0.79375000000000007
$ws.Range("B12").NumberFormat = "h:mm"

$ws.Range("C12").Value = 0.81736111111111109
$ws.Range("C12").NumberFormat = "h:mm"

# --- Row heights (auto-sized by Excel to fit the wrapped topic text) -------
$ws.Rows("12").RowHeight = 60
$ws.Rows("13").RowHeight = 37.5
$ws.Rows("14").RowHeight = 41.25
$ws.Rows("15").RowHeight = 28.5
$ws.Rows("16").RowHeight = 33
$ws.Rows("17").RowHeight = 36.75
$ws.Rows("18").RowHeight = 25.5

# --- View state: zoomed to 85%, scrolled down, selection on G17 ------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("G17").Select() | Out-Null
